$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.762.88'
$ws.Range("E2").Value = '  -3.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.908.12'
$ws.Range("E3").Value = '  -4.08%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.24'
$ws.Range("E5").Value = '  -1.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.27'
$ws.Range("E6").Value = '  -6.07%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -2.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.908.93'
$ws.Range("E9").Value = '  -3.94%  '
$ws.Range("E10").Value = '  -3.72%  '
$ws.Range("E11").Value = '  -5.42%  '
$ws.Range("E12").Value = '  -4.31%  '
$ws.Range("E13").Value = '  -3.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.41'
$ws.Range("E14").Value = '  -6.28%  '
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.387.53'
$ws.Range("E16").Value = '  -4.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.728.78'
$ws.Range("E17").Value = '  -3.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.67'
$ws.Range("E18").Value = '  -6.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.906.76'
$ws.Range("E19").Value = '  -4.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '427.70'
$ws.Range("E20").Value = '  -5.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.51'
$ws.Range("E21").Value = '  -5.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.681'
$ws.Range("E22").Value = '  -2.49%  '
$ws.Range("E23").Value = '  -6.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.76'
$ws.Range("E24").Value = '  -2.91%  '
$ws.Range("E25").Value = '  -5.57%  '
$ws.Range("E26").Value = '  -5.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.85'
$ws.Range("E27").Value = '  -4.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -4.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.19'
$ws.Range("E31").Value = '  -4.08%  '
$ws.Range("E32").Value = '  -3.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.37'
$ws.Range("E33").Value = '  -4.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.107'
$ws.Range("E34").Value = '  -3.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0856'
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("E36").Value = '  -2.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.57'
$ws.Range("E37").Value = '  -5.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.01'
$ws.Range("E38").Value = '  -4.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.36'
$ws.Range("E39").Value = '  -2.09%  '
$ws.Range("E40").Value = '  -5.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.98'
$ws.Range("E41").Value = '  -6.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.57'
$ws.Range("E42").Value = '  -5.76%  '
$ws.Range("E43").Value = '  -4.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.24'
$ws.Range("E44").Value = '  -7.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0349'
$ws.Range("E45").Value = '  -3.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '373.45'
$ws.Range("E46").Value = '  -5.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.693.04'
$ws.Range("E47").Value = '  -1.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.22'
$ws.Range("E48").Value = '  -0.80%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.00'
$ws.Range("E50").Value = '  -6.46%  '
$ws.Range("E51").Value = '  -3.01%  '
